{"js": "// Load all body paragraphs with their text so we can locate the ones we need\n// to edit by content (more robust than relying on fixed indices).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Step 1 -----------------------------------------------------------\n// Find the paragraph that ends the intro (\"...pgAdmin/postgres. \") and\n// append a new run with the \"Luis has been tasked...\" sentence, exactly as\n// a new trailing run (it keeps the paragraph's existing run(s) untouched).\nconst introText = \"pgAdmin/postgres\";\nconst luisSentence =\n  \"Luis has been tasked with creating the tables in Postgres, using the csv file headers and see which columns he finds most relevant.\";\n\nlet introParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(introText) !== -1) {\n    introParaIndex = i;\n    break;\n  }\n}\nif (introParaIndex === -1) {\n  throw new Error(\"Could not locate the intro paragraph containing 'pgAdmin/postgres'.\");\n}\nconst introPara = paragraphs.items[introParaIndex];\nintroPara.insertText(luisSentence, \"End\");\nawait context.sync();\n\n// --- Step 2 -----------------------------------------------------------\n// Insert an extra empty paragraph right after the (first) blank paragraph\n// that follows the intro paragraph, so that paragraph is now followed by\n// two blank paragraphs instead of one.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst blankAfterIntro = paragraphs.items[introParaIndex + 1];\nblankAfterIntro.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// --- Step 3 -----------------------------------------------------------\n// Locate the paragraph that still holds the original \"Luis has been\n// tasked...\" sentence (now located further down, right after\n// \"Extraction:\") and replace its content with the new Jupyter Notebook /\n// pandas DataFrames text. We rebuild the paragraph via OOXML so we can\n// reproduce the same run layout as Word would (three runs, with the\n// \"Jupyter\" run wrapped in proofErr spell-check markers, matching how the\n// rest of the document already marks that word elsewhere).\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet luisParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === luisSentence) {\n    luisParaIndex = i;\n    break;\n  }\n}\nif (luisParaIndex === -1) {\n  throw new Error(\"Could not locate the 'Luis has been tasked...' paragraph to replace.\");\n}\nconst luisPara = paragraphs.items[luisParaIndex];\nconst luisRange = luisPara.getRange(\"Whole\");\n\nconst replacementOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">The dataset was pulled from Kaggle, as mentioned before, and we will be using </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>Jupyter</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> Notebook is extract those dataset .csv files and put them into pandas DataFrames.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nluisRange.insertOoxml(replacementOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$luisSentence = \"Luis has been tasked with creating the tables in Postgres, using the csv file headers and see which columns he finds most relevant.\"\n\n# --- Step 1 -----------------------------------------------------------\n# Find the paragraph that ends the intro (\"...pgAdmin/postgres. \") and\n# append a new run with the \"Luis has been tasked...\" sentence right after\n# the existing text (leaving the existing run(s) untouched).\n$paras = $d.Paragraphs\n$introIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text -like \"*pgAdmin/postgres*\") {\n        $introIndex = $i\n        break\n    }\n}\nif ($introIndex -eq -1) {\n    throw \"Could not locate the intro paragraph containing 'pgAdmin/postgres'.\"\n}\n\n$introRange = $paras.Item($introIndex).Range\n# Collapse to just before the paragraph mark so the new text becomes part\n# of the same paragraph (as a new trailing run) instead of creating a new one.\n$introRange.SetRange($introRange.End - 1, $introRange.End - 1)\n$introRange.InsertAfter($luisSentence)\n\n# --- Step 2 -----------------------------------------------------------\n# Insert an extra empty paragraph right after the (first) blank paragraph\n# that follows the intro paragraph, so it is now followed by two blank\n# paragraphs instead of one.\n$paras = $d.Paragraphs\n$blankAfterIntro = $paras.Item($introIndex + 1)\n$blankAfterIntro.Range.InsertParagraphAfter()\n\n# --- Step 3 -----------------------------------------------------------\n# Locate the paragraph that still holds the original \"Luis has been\n# tasked...\" sentence (now further down, right after \"Extraction:\") and\n# replace its content with the new Jupyter Notebook / pandas DataFrames\n# text. Rebuild the paragraph via InsertXML so the run layout matches\n# Word's own output (three runs, with \"Jupyter\" wrapped in proofErr\n# spell-check markers, matching how the rest of the document already\n# marks that word elsewhere).\n$paras = $d.Paragraphs\n$luisIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text.Trim() -eq $luisSentence) {\n        $luisIndex = $i\n        break\n    }\n}\nif ($luisIndex -eq -1) {\n    throw \"Could not locate the 'Luis has been tasked...' paragraph to replace.\"\n}\n\n$luisRange = $paras.Item($luisIndex).Range\n# Exclude the trailing paragraph mark, then delete the old text before\n# inserting the new XML fragment so the content is replaced (not appended).\n$luisRange.SetRange($luisRange.Start, $luisRange.End - 1)\n$luisRange.Delete()\n\n$replacementXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">The dataset was pulled from Kaggle, as mentioned before, and we will be using </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Jupyter</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> Notebook is extract those dataset .csv files and put them into pandas DataFrames.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$luisRange.InsertXML($replacementXml)\n\n$d.Save()\n"}
